$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Replace the "Mifos style" value in B17 with the new scenario description,
# and give it a left/top-aligned, non-wrapping style (new cellXf).
$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160
$cell.WrapText = $false

# Leave the edited cell selected, matching the new active selection.
$cell.Select()
